# Updated cryptos list on Wed Oct  9 20:30:33 UTC 2024 with GitHub Actions
#
# This script applies the latest scraped price/volume values to the
# cryptocurrency table on Sheet1 (columns: A=rank index, B=Coin, C=Link,
# D=Price, E=Volume(1h)). Rows 50/51 (InjectiveProtocol / Hedera) also swap
# rank order, so their B/C/D/E cells are fully replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '60.885.33' },
    @{ Cell = 'E2'; Value = '  -2.16%  ' },
    @{ Cell = 'D3'; Value = '2.425.26' },
    @{ Cell = 'E3'; Value = '  -0.72%  ' },
    @{ Cell = 'D4'; Value = '0.998' },
    @{ Cell = 'D5'; Value = '569.65' },
    @{ Cell = 'E5'; Value = '  -2.30%  ' },
    @{ Cell = 'D6'; Value = '139.46' },
    @{ Cell = 'E6'; Value = '  -3.10%  ' },
    @{ Cell = 'E7'; Value = '  +0.19%  ' },
    @{ Cell = 'D8'; Value = '0.527' },
    @{ Cell = 'E8'; Value = '  -0.78%  ' },
    @{ Cell = 'D9'; Value = '2.415.15' },
    @{ Cell = 'E9'; Value = '  -1.02%  ' },
    @{ Cell = 'D10'; Value = '0.108' },
    @{ Cell = 'E10'; Value = '  +0.40%  ' },
    @{ Cell = 'E11'; Value = '  +0.18%  ' },
    @{ Cell = 'D12'; Value = '5.06' },
    @{ Cell = 'E12'; Value = '  -2.93%  ' },
    @{ Cell = 'D13'; Value = '0.336' },
    @{ Cell = 'E13'; Value = '  -2.46%  ' },
    @{ Cell = 'D14'; Value = '25.92' },
    @{ Cell = 'E14'; Value = '  -2.02%  ' },
    @{ Cell = 'D15'; Value = '0.0000169' },
    @{ Cell = 'E15'; Value = '  -2.13%  ' },
    @{ Cell = 'D16'; Value = '2.809.92' },
    @{ Cell = 'E16'; Value = '  -1.90%  ' },
    @{ Cell = 'D17'; Value = '61.022.29' },
    @{ Cell = 'E17'; Value = '  -1.74%  ' },
    @{ Cell = 'D18'; Value = '2.407.11' },
    @{ Cell = 'E18'; Value = '  -1.07%  ' },
    @{ Cell = 'D19'; Value = '10.49' },
    @{ Cell = 'E19'; Value = '  -3.75%  ' },
    @{ Cell = 'D20'; Value = '7.24' },
    @{ Cell = 'E20'; Value = '  +1.34%  ' },
    @{ Cell = 'D21'; Value = '321.46' },
    @{ Cell = 'E21'; Value = '  -2.72%  ' },
    @{ Cell = 'D22'; Value = '4.01' },
    @{ Cell = 'E22'; Value = '  -2.33%  ' },
    @{ Cell = 'E23'; Value = '  +2.11%  ' },
    @{ Cell = 'E24'; Value = '  +0.09%  ' },
    @{ Cell = 'D25'; Value = '1.86' },
    @{ Cell = 'E25'; Value = '  -5.96%  ' },
    @{ Cell = 'D26'; Value = '64.77' },
    @{ Cell = 'E26'; Value = '  -1.53%  ' },
    @{ Cell = 'D27'; Value = '8.68' },
    @{ Cell = 'E27'; Value = '  -7.53%  ' },
    @{ Cell = 'D28'; Value = '571.57' },
    @{ Cell = 'E28'; Value = '  -7.82%  ' },
    @{ Cell = 'E29'; Value = '  +0.10%  ' },
    @{ Cell = 'D30'; Value = '0.0₃0902' },
    @{ Cell = 'E30'; Value = '  -5.74%  ' },
    @{ Cell = 'D31'; Value = '7.80' },
    @{ Cell = 'E31'; Value = '  -2.46%  ' },
    @{ Cell = 'D32'; Value = '1.33' },
    @{ Cell = 'E32'; Value = '  -7.12%  ' },
    @{ Cell = 'E33'; Value = '  -3.00%  ' },
    @{ Cell = 'D34'; Value = '0.131' },
    @{ Cell = 'E34'; Value = '  -6.92%  ' },
    @{ Cell = 'E35'; Value = '  +0.27%  ' },
    @{ Cell = 'D36'; Value = '4.56' },
    @{ Cell = 'E36'; Value = '  -7.36%  ' },
    @{ Cell = 'E37'; Value = '  -3.61%  ' },
    @{ Cell = 'D38'; Value = '148.39' },
    @{ Cell = 'E38'; Value = '  -2.08%  ' },
    @{ Cell = 'D39'; Value = '1.36' },
    @{ Cell = 'E39'; Value = '  -4.71%  ' },
    @{ Cell = 'D40'; Value = '18.13' },
    @{ Cell = 'E40'; Value = '  -1.09%  ' },
    @{ Cell = 'D41'; Value = '5.03' },
    @{ Cell = 'E41'; Value = '  -4.38%  ' },
    @{ Cell = 'E42'; Value = '  +0.00%  ' },
    @{ Cell = 'E43'; Value = '  -1.71%  ' },
    @{ Cell = 'D44'; Value = '1.64' },
    @{ Cell = 'E44'; Value = '  -6.59%  ' },
    @{ Cell = 'D45'; Value = '2.32' },
    @{ Cell = 'E45'; Value = '  -6.52%  ' },
    @{ Cell = 'D46'; Value = '0.0₆0277' },
    @{ Cell = 'E46'; Value = '  +17.64%  ' },
    @{ Cell = 'D47'; Value = '140.05' },
    @{ Cell = 'E47'; Value = '  -2.29%  ' },
    @{ Cell = 'D48'; Value = '3.48' },
    @{ Cell = 'E48'; Value = '  -4.15%  ' },
    @{ Cell = 'D49'; Value = '0.592' },
    @{ Cell = 'E49'; Value = '  -0.89%  ' },
    @{ Cell = 'B50'; Value = 'Hedera' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D50'; Value = '0.0502' },
    @{ Cell = 'E50'; Value = '  -4.37%  ' },
    @{ Cell = 'B51'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D51'; Value = '19.18' },
    @{ Cell = 'E51'; Value = '  -1.64%  ' }
)

foreach ($change in $changes) {
    $cellRef = $change.Cell
    $val = $change.Value

    # The Price column holds plain-looking decimal numbers as TEXT
    # (e.g. "0.998", "7.80") so trailing zeros and exact formatting are
    # preserved. A bare assignment of a numeric-looking string lets Excel
    # auto-convert the cell to a Number (dropping trailing zeros / type),
    # so such values are written with a leading apostrophe to force text,
    # matching the original cell's text semantics.
    if ($val -match '^[0-9]+(\.[0-9]+)?$') {
        $ws.Range($cellRef).Value = "'" + $val
    } else {
        $ws.Range($cellRef).Value = $val
    }
}
